$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1568.9333
$ws.Range("I125").Value = 1325
$ws.Range("J125").Value = 1657.6364
$ws.Range("K125").Value = 11925
$ws.Range("L125").Value = 14918.7276
$ws.Range("M125").Value = -9465
$ws.Range("N125").Value = -19838.7276

$ws.Range("H131").Value = 2621.65
$ws.Range("I131").Value = 2074.7778
$ws.Range("J131").Value = 3069.0908
$ws.Range("K131").Value = 6224.3334
$ws.Range("L131").Value = 9207.2724
$ws.Range("M131").Value = -1184.3334
$ws.Range("N131").Value = -19287.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 490.5
$ws.Range("I22").Value = 517.7143
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 517.7143
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -344.7143
$ws.Range("N22").Value = -646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1677.6666
$ws.Range("I16").Value = 1570.2222
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1570.2222
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1283.2222
$ws.Range("N16").Value = -2574

$ws.Range("H86").Value = 3186.7778
$ws.Range("I86").Value = 3042
$ws.Range("J86").Value = 3563.2
$ws.Range("K86").Value = 3042
$ws.Range("L86").Value = 3563.2
$ws.Range("M86").Value = -1919
$ws.Range("N86").Value = -5809.2

$ws.Range("H89").Value = 3186.7778
$ws.Range("I89").Value = 3042
$ws.Range("J89").Value = 3563.2
$ws.Range("K89").Value = 15210
$ws.Range("L89").Value = 17816
$ws.Range("M89").Value = -9594
$ws.Range("N89").Value = -29048

$ws.Range("H99").Value = 3382.3157
$ws.Range("I99").Value = 3672.3333
$ws.Range("J99").Value = 3248.4614
$ws.Range("K99").Value = 3672.3333
$ws.Range("L99").Value = 3248.4614
$ws.Range("M99").Value = -2174.3333
$ws.Range("N99").Value = -6244.4614

$ws.Range("H113").Value = 1677.6666
$ws.Range("I113").Value = 1570.2222
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1570.2222
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 599.7778000000001
$ws.Range("N113").Value = -6340

$ws.Range("H126").Value = 3382.3157
$ws.Range("I126").Value = 3672.3333
$ws.Range("J126").Value = 3248.4614
$ws.Range("K126").Value = 11016.9999
$ws.Range("L126").Value = 9745.3842
$ws.Range("M126").Value = -8546.999899999999
$ws.Range("N126").Value = -14685.3842

$ws.Range("H141").Value = 3198.6667
$ws.Range("J141").Value = 3198.6667
$ws.Range("L141").Value = 3198.6667
$ws.Range("N141").Value = -13558.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 728.5714
$ws.Range("J32").Value = 728.5714
$ws.Range("L32").Value = 2185.7142
$ws.Range("N32").Value = -2751.7142

$ws.Range("H107").Value = 8996.392
$ws.Range("I107").Value = 6577.875
$ws.Range("J107").Value = 14524.429
$ws.Range("K107").Value = 19733.625
$ws.Range("L107").Value = 43573.287
$ws.Range("M107").Value = -17813.625
$ws.Range("N107").Value = -47413.287

$ws.Range("H113").Value = 5759.8
$ws.Range("I113").Value = 9621.091
$ws.Range("J113").Value = 1040.4445
$ws.Range("K113").Value = 28863.273
$ws.Range("L113").Value = 3121.3335
$ws.Range("M113").Value = -26693.273
$ws.Range("N113").Value = -7461.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 753056
$ws.Range("I2").Value = 1204854.8
$ws.Range("J2").Value = 58.166668
$ws.Range("K2").Value = 1204854.8
$ws.Range("L2").Value = 58.166668
$ws.Range("M2").Value = -1204741.8
$ws.Range("N2").Value = -284.166668

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H62").Value = 21983.334

$ws.Range("H65").Value = 21983.334

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H122").Value = 1677.6666
$ws.Range("I122").Value = 1699.8823
$ws.Range("K122").Value = 5099.6469
$ws.Range("M122").Value = -2649.6469

$ws.Range("H132").Value = 22225818
$ws.Range("I132").Value = 35717650
$ws.Range("J132").Value = 3978.5293
$ws.Range("K132").Value = 107152950
$ws.Range("L132").Value = 11935.5879
$ws.Range("M132").Value = -107150420
$ws.Range("N132").Value = -16995.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 960
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1790

$ws.Range("H27").Value = 960
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1414

$ws.Range("H40").Value = 3535.12
$ws.Range("I40").Value = 2743.9375
$ws.Range("J40").Value = 4941.6665
$ws.Range("K40").Value = 2743.9375
$ws.Range("L40").Value = 4941.6665
$ws.Range("M40").Value = -2607.9375
$ws.Range("N40").Value = -5213.6665

$ws.Range("H100").Value = 2263.111
$ws.Range("I100").Value = 2296
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2296
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1755
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3598.125
$ws.Range("I81").Value = 966.9231
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 1933.8462
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -872.8462
$ws.Range("N81").Value = -32122

$ws.Range("H84").Value = 3598.125
$ws.Range("I84").Value = 966.9231
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 9669.231
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -4365.231
$ws.Range("N84").Value = -160608

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H107").Value = 7693454.5
$ws.Range("I107").Value = 614
$ws.Range("J107").Value = 16668435
$ws.Range("K107").Value = 1842
$ws.Range("L107").Value = 50005305
$ws.Range("M107").Value = 78
$ws.Range("N107").Value = -50009145

$ws.Range("H113").Value = 1160.5
$ws.Range("I113").Value = 1222
$ws.Range("J113").Value = 1099
$ws.Range("K113").Value = 3666
$ws.Range("L113").Value = 3297
$ws.Range("M113").Value = -1496
$ws.Range("N113").Value = -7637

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H126").Value = 1084.8298
$ws.Range("I126").Value = 776.5714
$ws.Range("J126").Value = 1983.9166
$ws.Range("K126").Value = 2329.7142
$ws.Range("L126").Value = 5951.7498
$ws.Range("M126").Value = 140.2857999999997
$ws.Range("N126").Value = -10891.7498
